# Weekly crime-data refresh: NYPD 9th Precinct CompStat report.
# - New Police Commissioner byline.
# - Report header advances one week (Volume/Number + date range).
# - Per-category Week-to-Date / 28-Day / Year-to-Date / multi-year figures
#   refreshed with newly collected counts (and their derived % changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text / shared-string updates ---
$ws.Range("M6").Value = "Edward A. Caban"
$ws.Range("A8").Value = "Volume 30   Number  27"
$ws.Range("C9").Value = "Report Covering the Week  7/3/2023  Through  7/9/2023"

# --- Cells whose type/style changes (number <-> text placeholder) ---
# Use Copy-to-destination so the destination inherits the exact cellXf (style index)
# used elsewhere in the sheet for that data type, matching the target workbook.
$ws.Range("D26").Copy($ws.Range("C20"))   # -> text "0" style (s=14)
$ws.Range("D26").Copy($ws.Range("D20"))   # -> text "0" style (s=14)
$ws.Range("E14").Copy($ws.Range("E20"))   # -> text "***.*" style (s=14)
$ws.Range("D26").Copy($ws.Range("C27"))   # -> text "0" style (s=14)
$ws.Range("J15").Copy($ws.Range("F30"))   # -> numeric style (s=15)
$ws.Range("J15").Copy($ws.Range("I30"))   # -> numeric style (s=15)
$ws.Range("F30").Value = 2
$ws.Range("I30").Value = 2

# --- Numeric value updates (style/type unchanged) ---
$ws.Range("C15").Value = 2
$ws.Range("F15").Value = 4
$ws.Range("I15").Value = 12
$ws.Range("K15").Value = 9.090909090909
$ws.Range("M15").Value = 140
$ws.Range("N15").Value = -52
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = -81.818181818181
$ws.Range("F16").Value = 12
$ws.Range("G16").Value = 31
$ws.Range("H16").Value = -61.290322580645
$ws.Range("I16").Value = 91
$ws.Range("J16").Value = 142
$ws.Range("K16").Value = -35.915492957746
$ws.Range("L16").Value = 30
$ws.Range("M16").Value = 5.813953488372
$ws.Range("N16").Value = -79.458239277652
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 66.666666666666
$ws.Range("I17").Value = 131
$ws.Range("J17").Value = 98
$ws.Range("K17").Value = 33.673469387755
$ws.Range("L17").Value = 50.574712643678
$ws.Range("M17").Value = 57.831325301204
$ws.Range("N17").Value = -57.467532467532
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -42.857142857142
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 27
$ws.Range("H18").Value = -29.629629629629
$ws.Range("I18").Value = 145
$ws.Range("J18").Value = 202
$ws.Range("K18").Value = -28.217821782178
$ws.Range("L18").Value = -9.375
$ws.Range("M18").Value = 22.881355932203
$ws.Range("N18").Value = -66.435185185185
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 68
$ws.Range("G19").Value = 79
$ws.Range("H19").Value = -13.924050632911
$ws.Range("I19").Value = 517
$ws.Range("J19").Value = 501
$ws.Range("K19").Value = 3.193612774451
$ws.Range("L19").Value = 61.059190031152
$ws.Range("M19").Value = 32.225063938618
$ws.Range("N19").Value = -31.794195250659
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 150
$ws.Range("L20").Value = -46.511627906976
$ws.Range("M20").Value = -8
$ws.Range("N20").Value = -91.481481481481
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 44
$ws.Range("E21").Value = -43.181818181818
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 151
$ws.Range("H21").Value = -15.231788079470
$ws.Range("I21").Value = 920
$ws.Range("J21").Value = 982
$ws.Range("K21").Value = -6.313645621181
$ws.Range("L21").Value = 32.374100719424
$ws.Range("M21").Value = 29.577464788732
$ws.Range("N21").Value = -58.946898705934
$ws.Range("G22").Value = 3
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -20
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = -33.333333333333
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 11.111111111111
$ws.Range("I23").Value = 72
$ws.Range("J23").Value = 74
$ws.Range("K23").Value = -2.702702702702
$ws.Range("L23").Value = -30.769230769230
$ws.Range("M23").Value = 10.769230769230
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 41
$ws.Range("E24").Value = -43.902439024390
$ws.Range("F24").Value = 135
$ws.Range("G24").Value = 179
$ws.Range("H24").Value = -24.581005586592
$ws.Range("I24").Value = 762
$ws.Range("J24").Value = 1171
$ws.Range("K24").Value = -34.927412467976
$ws.Range("L24").Value = 60.421052631578
$ws.Range("M24").Value = -11.805555555555
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 37
$ws.Range("G25").Value = 42
$ws.Range("H25").Value = -11.904761904761
$ws.Range("I25").Value = 247
$ws.Range("J25").Value = 245
$ws.Range("K25").Value = 0.816326530612
$ws.Range("L25").Value = 40.340909090909
$ws.Range("M25").Value = 4.661016949152
$ws.Range("C26").Value = 3
$ws.Range("F26").Value = 5
$ws.Range("I26").Value = 21
$ws.Range("K26").Value = 10.526315789473
$ws.Range("L26").Value = 16.666666666666
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = -100
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = 16.666666666666
$ws.Range("J27").Value = 49
$ws.Range("K27").Value = -40.816326530612
$ws.Range("L27").Value = -6.451612903225
$ws.Range("N28").Value = -84.210526315789
$ws.Range("N29").Value = -80
$ws.Range("K30").Value = -71.428571428571
$ws.Range("L30").Value = -66.666666666666
